$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the dataset (original row numbers,
# deleted bottom-first so the earlier row index stays valid):
#   row 28 -> "SC 92"
#   row 26 -> "RM 232"
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# Cell-level value corrections (row numbers below are the ones in the
# resulting, already-shifted sheet).
$ws.Range("E5").ClearContents()
$ws.Range("E11").Value = -7.9
$ws.Range("C19").Value = 13.2
$ws.Range("E19").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("C23").Value = 12.2
$ws.Range("E25").Value = -7.1
$ws.Range("C27").ClearContents()
$ws.Range("F28").Value = 17.44
$ws.Range("E29").ClearContents()
$ws.Range("F32").ClearContents()
$ws.Range("C33").Value = 10.4
